$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "67.533.86"
$ws.Range("E2").Value = "  -2.86%  "
Set-TextValue $ws.Range("D3") "3.506.66"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "606.32"
$ws.Range("E5").Value = "  -1.66%  "
Set-TextValue $ws.Range("D6") "149.81"
$ws.Range("E6").Value = "  -5.87%  "
Set-TextValue $ws.Range("D7") "3.508.25"
$ws.Range("E7").Value = "  -4.36%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("E10").Value = "  -3.52%  "
Set-TextValue $ws.Range("D11") "6.97"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("E13").Value = "  -3.86%  "
Set-TextValue $ws.Range("D14") "4.098.03"
$ws.Range("E14").Value = "  -4.49%  "
Set-TextValue $ws.Range("D15") "31.57"
$ws.Range("E15").Value = "  -2.32%  "
Set-TextValue $ws.Range("D16") "3.502.48"
$ws.Range("E16").Value = "  -4.79%  "
Set-TextValue $ws.Range("D17") "67.467.63"
$ws.Range("E17").Value = "  -3.08%  "
$ws.Range("E18").Value = "  -0.57%  "
Set-TextValue $ws.Range("D19") "6.37"
$ws.Range("E19").Value = "  -1.92%  "
Set-TextValue $ws.Range("D20") "15.01"
$ws.Range("E20").Value = "  -5.39%  "
Set-TextValue $ws.Range("D21") "447.25"
$ws.Range("E21").Value = "  -4.61%  "
Set-TextValue $ws.Range("D22") "8.98"
$ws.Range("E22").Value = "  -12.61%  "
Set-TextValue $ws.Range("D23") "0.620"
$ws.Range("E23").Value = "  -4.14%  "
Set-TextValue $ws.Range("D24") "77.37"
$ws.Range("E24").Value = "  -2.58%  "
Set-TextValue $ws.Range("D25") "0.0000129"
$ws.Range("E25").Value = "  +6.17%  "
$ws.Range("E26").Value = "  +0.11%  "
Set-TextValue $ws.Range("D27") "3.645.05"
$ws.Range("E27").Value = "  -4.50%  "
Set-TextValue $ws.Range("D28") "10.19"
$ws.Range("E28").Value = "  -7.60%  "
Set-TextValue $ws.Range("D29") "8.26"
$ws.Range("E29").Value = "  -5.03%  "
Set-TextValue $ws.Range("D30") "2.48"
$ws.Range("E30").Value = "  -4.18%  "
Set-TextValue $ws.Range("D31") "1.55"
$ws.Range("E31").Value = "  -6.53%  "
$ws.Range("E32").Value = "  -0.03%  "
Set-TextValue $ws.Range("D33") "0.165"
$ws.Range("E33").Value = "  +1.29%  "
Set-TextValue $ws.Range("D34") "25.72"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("E35").Value = "  -3.28%  "
Set-TextValue $ws.Range("D36") "1.85"
$ws.Range("E36").Value = "  -6.23%  "
Set-TextValue $ws.Range("D37") "3.493.73"
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -0.10%  "
Set-TextValue $ws.Range("D41") "2.20"
$ws.Range("E41").Value = "  -0.56%  "
Set-TextValue $ws.Range("D42") "173.07"
$ws.Range("E42").Value = "  -2.84%  "
Set-TextValue $ws.Range("D43") "0.0876"
$ws.Range("E43").Value = "  -1.47%  "
Set-TextValue $ws.Range("D44") "5.42"
$ws.Range("E44").Value = "  -5.88%  "
Set-TextValue $ws.Range("D45") "0.882"
$ws.Range("E45").Value = "  -4.72%  "
Set-TextValue $ws.Range("D46") "45.45"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "27.83"
$ws.Range("E47").Value = "  -4.25%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D48") "1.28"
$ws.Range("E48").Value = "  +6.51%  "
Set-TextValue $ws.Range("D49") "2.56"
$ws.Range("E49").Value = "  -4.75%  "
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("E51").Value = "  -2.91%  "
